# Rename "center of mass" columns to "centroid weighted" and add new
# geometric "centroid" columns on the FieldIlluminationKeyValues sheet.
#
# Before: A=channel, B=nb_pixels_center, C=center_of_mass_x, D=center_of_mass_y,
#         E=max_intensity, ...
# After:  A=channel, B=nb_pixels_center, C=centroid_weighted_y, D=centroid_weighted_x,
#         E=centroid_y, F=centroid_x, G=max_intensity, ...

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FieldIlluminationKeyValues")

# Insert two new blank columns right before the old "max_intensity" column
# (originally column E), pushing it (and everything after it) two columns
# to the right, while leaving C/D (center_of_mass_x/y) in place so they can
# be renamed below.
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).Insert()

# Rename the former center-of-mass columns to the new "centroid_weighted"
# naming (note the x/y swap, matching the upstream rename).
$ws.Range("C1").Value = "centroid_weighted_y"
$ws.Range("D1").Value = "centroid_weighted_x"

# Populate the two newly inserted columns with the new geometric centroid
# fields.
$ws.Range("E1").Value = "centroid_y"
$ws.Range("F1").Value = "centroid_x"
